$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.757.87'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '3.456.15'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.614'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.46%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '3.450.40'
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("E12").Value = '  -1.81%  '
$ws.Range("D13").Value = '4.053.56'
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("E15").Value = '  -1.09%  '
$ws.Range("D16").Value = '67.729.51'
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000176'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.11%  '
$ws.Range("D18").Value = '3.456.37'
$ws.Range("E18").Value = '  -1.39%  '
$ws.Range("E19").Value = '  -3.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.50%  '
$ws.Range("E22").Value = '  -2.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.82'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("E25").Value = '  -1.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000120'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.43%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.35'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.29%  '
$ws.Range("E29").Value = '  -1.75%  '
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("E31").Value = '  -3.37%  '
$ws.Range("E32").Value = '  -1.28%  '
$ws.Range("E33").Value = '  -5.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.92%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.33'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.86%  '
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.56'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '162.16'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("E40").Value = '  +1.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.87'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.26%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.62'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.99%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.68'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0716'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.12'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.21%  '
$ws.Range("D47").Value = '2.722.80'
$ws.Range("E47").Value = '  -4.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '41.25'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.75%  '
$ws.Range("E49").Value = '  -2.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '327.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.21%  '
$ws.Range("E51").Value = '  -4.37%  '
